{"js": "// Remove the sentence \"tenant nao esta furando accont e vice versa?\" together\n// with the manual line break that followed it, so the next sentence (\"tem\n// algum dom\u00ednio \"furando\" o outro, por exemplo ...\") continues right after\n// \"...migrations? \" in the same paragraph instead of starting on a new line.\n//\n// Word represents a manual line break (<w:br/>) as the vertical-tab\n// character (\"\\u000b\") in Range/Body text, so it can be matched by a normal\n// text search.\nconst body = context.document.body;\n\n// Include the leading \"? \" in the search/replace so the match starts at a\n// run boundary; this lets the host clean up the spell-check (proofErr)\n// marker that used to wrap the now-deleted \"tenant\" without leaving an\n// orphaned tag behind. We simply put the \"? \" back at the start of the\n// replacement text, so the visible content is unchanged.\nconst results = body.search(\"? tenant nao esta furando accont e vice versa?\\u000b\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"? \", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Remove the sentence \"tenant nao esta furando accont e vice versa?\" together\n# with the manual line break that followed it, so the next sentence (\"tem\n# algum dom\u00ednio \"furando\" o outro, por exemplo ...\") continues right after\n# \"...migrations? \" in the same paragraph instead of starting on a new line.\n#\n# Word represents a manual line break (<w:br/>) as Chr(11) (vertical tab) in\n# Range.Text, so it can be matched directly by Find.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n# Include the leading \"? \" in the Find/Replace text so the match starts at a\n# run boundary; this lets Word clean up the spell-check (proofErr) marker\n# that used to wrap the now-deleted \"tenant\" without leaving an orphaned tag\n# behind. We put \"? \" back at the start of the replacement text, so the\n# visible content is unchanged.\n$findText = \"? tenant nao esta furando accont e vice versa?\" + [char]11\n\n$found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, \"? \", 2)\n\nif ($found) {\n    Write-Output \"Replaced leak sentence + manual line break.\"\n} else {\n    Write-Output \"Target text not found; document left unchanged.\"\n}\n"}
